$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(10).EntireColumn.AutoFit()
Write-Host "J width after autofit:" $ws.Columns.Item(10).ColumnWidth
